# Se agrega en el streamlit la info de competencias
#
# The workbook tracks weekly competitor pricing ("ANIO"/"SEMNUMERO" = year /
# week number). Rows 107:148 on Hoja1 still say ANIO=2024 for weeks 2-43
# even though row 106 (week 1) already rolled over to 2025 - fix the rest
# of the year column to read 2025 as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Fix the ANIO (year) column for weeks 2-43: 2024 -> 2025 ---
$ws.Range("C107:C148").Value = 2025

# --- Column D ("SEMNUMERO") now holds 2-digit week numbers; widen it to fit ---
$ws.Columns("D:D").AutoFit()

# --- Leave the view scrolled down to / focused on the rows that changed ---
$ws.Range("C106:C148").Select()
